$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pcs")

# Insert a "total" row after each demographic sub-group, summing the
# percentage columns (G = before, H = after) for that group, and bold
# the new total cells to match the rest of the workbook's group-header
# styling.

# Group 1 (rows 3-7)
$ws.Rows.Item(8).Insert()
$ws.Range("G8").Formula = "=SUM(G3:G7)"
$ws.Range("H8").Formula = "=SUM(H3:H7)"
$ws.Range("G8:H8").Font.Bold = $true

# Group 2 (rows 9-11 after previous insert)
$ws.Rows.Item(12).Insert()
$ws.Range("G12").Formula = "=SUM(G9:G11)"
$ws.Range("H12").Formula = "=SUM(H9:H11)"
$ws.Range("G12:H12").Font.Bold = $true

# Group 3 (rows 13-15 after previous inserts)
$ws.Rows.Item(16).Insert()
$ws.Range("G16").Formula = "=SUM(G13:G15)"
$ws.Range("H16").Formula = "=SUM(H13:H15)"
$ws.Range("G16:H16").Font.Bold = $true

# Group 4 (rows 17-18 after previous inserts)
$ws.Rows.Item(19).Insert()
$ws.Range("G19").Formula = "=SUM(G17:G18)"
$ws.Range("H19").Formula = "=SUM(H17:H18)"
$ws.Range("G19:H19").Font.Bold = $true

# Group 5 (rows 20-21 after previous inserts)
$ws.Rows.Item(22).Insert()
$ws.Range("G22").Formula = "=SUM(G20:G21)"
$ws.Range("H22").Formula = "=SUM(H20:H21)"
$ws.Range("G22:H22").Font.Bold = $true

# Group 6 (rows 23-27 after previous inserts)
$ws.Rows.Item(28).Insert()
$ws.Range("G28").Formula = "=SUM(G23:G27)"
$ws.Range("H28").Formula = "=SUM(H23:H27)"
$ws.Range("G28:H28").Font.Bold = $true
# Row insert copied the neighbouring I27 style down into I28 even though
# that column has no total; drop the stray formatted-but-empty cell.
$ws.Range("I28").Clear()

# Restore view state: zoomed to 90%, selection parked below the new data
$ws.Activate()
$excel.ActiveWindow.Zoom = 90
$ws.Range("G31").Select()
